$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells stay text (values contain multiple dots, e.g. thousands separators)
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "28.114.29"
$ws.Cells.Item(2, 5).Value = "  +0.11%  "

$ws.Cells.Item(3, 4).Value = "1.798.54"
$ws.Cells.Item(3, 5).Value = "  +2.44%  "

$ws.Cells.Item(4, 4).Value = "1.005"
$ws.Cells.Item(4, 5).Value = "  +0.19%  "

$ws.Cells.Item(5, 4).Value = "337.31"
$ws.Cells.Item(5, 5).Value = "  +0.05%  "

$ws.Cells.Item(6, 5).Value = "  +0.10%  "

$ws.Cells.Item(7, 4).Value = "0.4690"
$ws.Cells.Item(7, 5).Value = "  +24.30%  "

$ws.Cells.Item(8, 4).Value = "0.3712"
$ws.Cells.Item(8, 5).Value = "  +10.75%  "

$ws.Cells.Item(9, 4).Value = "45.28"
$ws.Cells.Item(9, 5).Value = "  -0.05%  "

$ws.Cells.Item(10, 4).Value = "0.07660"
$ws.Cells.Item(10, 5).Value = "  +6.51%  "

$ws.Cells.Item(11, 5).Value = "  +2.78%  "

$ws.Cells.Item(12, 4).Value = "22.61"
$ws.Cells.Item(12, 5).Value = "  +0.20%  "

$ws.Cells.Item(13, 4).Value = "1.004"
$ws.Cells.Item(13, 5).Value = "  +0.23%  "

$ws.Cells.Item(14, 4).Value = "6.352"
$ws.Cells.Item(14, 5).Value = "  +3.58%  "

$ws.Cells.Item(15, 4).Value = "7.385"
$ws.Cells.Item(15, 5).Value = "  +3.48%  "

$ws.Cells.Item(16, 4).Value = "1.797.38"
$ws.Cells.Item(16, 5).Value = "  +2.37%  "

$ws.Cells.Item(17, 4).Value = "0.00001094"
$ws.Cells.Item(17, 5).Value = "  +3.64%  "

$ws.Cells.Item(18, 4).Value = "0.06727"
$ws.Cells.Item(18, 5).Value = "  +2.16%  "

$ws.Cells.Item(19, 4).Value = "82.56"
$ws.Cells.Item(19, 5).Value = "  +2.71%  "

$ws.Cells.Item(21, 4).Value = "17.41"
$ws.Cells.Item(21, 5).Value = "  +3.10%  "

$ws.Cells.Item(22, 4).Value = "6.411"
$ws.Cells.Item(22, 5).Value = "  +2.85%  "

$ws.Cells.Item(23, 4).Value = "28.125.62"
$ws.Cells.Item(23, 5).Value = "  +0.16%  "

$ws.Cells.Item(24, 4).Value = "11.88"
$ws.Cells.Item(24, 5).Value = "  +2.04%  "

$ws.Cells.Item(25, 5).Value = "  +0.40%  "

$ws.Cells.Item(26, 4).Value = "20.84"
$ws.Cells.Item(26, 5).Value = "  +5.29%  "

$ws.Cells.Item(27, 4).Value = "2.386"
$ws.Cells.Item(27, 5).Value = "  +2.93%  "

$ws.Cells.Item(28, 4).Value = "151.58"
$ws.Cells.Item(28, 5).Value = "  -0.66%  "

$ws.Cells.Item(29, 4).Value = "2.003.56"
$ws.Cells.Item(29, 5).Value = "  +2.39%  "

$ws.Cells.Item(30, 4).Value = "133.80"
$ws.Cells.Item(30, 5).Value = "  +1.59%  "

$ws.Cells.Item(31, 4).Value = "1.260"
$ws.Cells.Item(31, 5).Value = "  +0.70%  "

$ws.Cells.Item(32, 5).Value = "  +0.60%  "

$ws.Cells.Item(33, 4).Value = "0.09705"
$ws.Cells.Item(33, 5).Value = "  +11.24%  "

$ws.Cells.Item(34, 4).Value = "5.922"
$ws.Cells.Item(34, 5).Value = "  +2.59%  "

$ws.Cells.Item(35, 4).Value = "0.02380"
$ws.Cells.Item(35, 5).Value = "  +2.16%  "

$ws.Cells.Item(36, 2).Value = "Algorand"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(36, 4).Value = "0.2228"
$ws.Cells.Item(36, 5).Value = "  +5.64%  "

$ws.Cells.Item(37, 2).Value = "Aptos"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(37, 4).Value = "12.19"
$ws.Cells.Item(37, 5).Value = "  -0.18%  "

$ws.Cells.Item(38, 4).Value = "0.06351"
$ws.Cells.Item(38, 5).Value = "  +2.76%  "

$ws.Cells.Item(39, 4).Value = "0.6718"
$ws.Cells.Item(39, 5).Value = "  +0.78%  "

$ws.Cells.Item(40, 4).Value = "5.264"
$ws.Cells.Item(40, 5).Value = "  +2.07%  "

$ws.Cells.Item(41, 2).Value = "WEMIXTOKEN"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(41, 4).Value = "1.503"
$ws.Cells.Item(41, 5).Value = "  +4.17%  "

$ws.Cells.Item(42, 2).Value = "TrustWalletToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(42, 4).Value = "1.238"
$ws.Cells.Item(42, 5).Value = "  +1.96%  "

$ws.Cells.Item(43, 4).Value = "8.071"
$ws.Cells.Item(43, 5).Value = "  +0.61%  "

$ws.Cells.Item(44, 5).Value = "  +3.77%  "

$ws.Cells.Item(45, 4).Value = "1.001"

$ws.Cells.Item(46, 4).Value = "0.6166"
$ws.Cells.Item(46, 5).Value = "  +2.03%  "

$ws.Cells.Item(47, 4).Value = "3.854"
$ws.Cells.Item(47, 5).Value = "  +0.44%  "

$ws.Cells.Item(48, 4).Value = "130.35"
$ws.Cells.Item(48, 5).Value = "  +1.42%  "

$ws.Cells.Item(49, 4).Value = "2.061"
$ws.Cells.Item(49, 5).Value = "  +2.34%  "

$ws.Cells.Item(50, 4).Value = "1.183"
$ws.Cells.Item(50, 5).Value = "  +0.70%  "

$ws.Cells.Item(51, 4).Value = "0.07133"
$ws.Cells.Item(51, 5).Value = "  -0.21%  "

# Restore default number formatting/style now that text values are set
$priceRange.NumberFormat = "General"
$priceRange.Style = "Normal"

